$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5500
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9830
$ws.Range("H19").Value = 1299.5
$ws.Range("J19").Value = 1399
$ws.Range("L19").Value = 1399
$ws.Range("N19").Value = -1749
$ws.Range("H86").Value = 95889.2
$ws.Range("I86").Value = 659.6
$ws.Range("K86").Value = 659.6
$ws.Range("M86").Value = 463.4
$ws.Range("H89").Value = 95889.2
$ws.Range("I89").Value = 659.6
$ws.Range("K89").Value = 3298
$ws.Range("M89").Value = 2318
$ws.Range("H98").Value = 2040.6666
$ws.Range("H100").Value = 499
$ws.Range("I100").Value = 499
$ws.Range("K100").Value = 499
$ws.Range("M100").Value = 42
$ws.Range("H107").Value = 347.25
$ws.Range("I107").Value = 347.25
$ws.Range("K107").Value = 347.25
$ws.Range("M107").Value = 1572.75
$ws.Range("H112").Value = 3528
$ws.Range("J112").Value = 3556.1052
$ws.Range("L112").Value = 10668.3156
$ws.Range("N112").Value = -12884.3156
$ws.Range("H122").Value = 2040.6666
$ws.Range("H137").Value = 1336.5
$ws.Range("I137").Value = 1114.2858
$ws.Range("K137").Value = 3342.8574
$ws.Range("M137").Value = -792.8574000000003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3503534.2
$ws.Range("I32").Value = 3337036
$ws.Range("K32").Value = 3337036
$ws.Range("M32").Value = -3336749
$ws.Range("H45").Value = 3972
$ws.Range("I45").Value = 3972
$ws.Range("K45").Value = 3972
$ws.Range("M45").Value = -3595
$ws.Range("H61").Value = 2188.6
$ws.Range("I61").Value = 2232.111
$ws.Range("K61").Value = 2232.111
$ws.Range("M61").Value = -2020.111
$ws.Range("H74").Value = 1327.2858
$ws.Range("I74").Value = 1327.2858
$ws.Range("K74").Value = 1327.2858
$ws.Range("M74").Value = -453.2858000000001
$ws.Range("H77").Value = 1327.2858
$ws.Range("I77").Value = 1327.2858
$ws.Range("K77").Value = 6636.429
$ws.Range("M77").Value = -2268.429
$ws.Range("H132").Value = 1048.25
$ws.Range("I132").Value = 1014.3333
$ws.Range("K132").Value = 3042.9999
$ws.Range("M132").Value = -512.9998999999998
$ws.Range("H136").Value = 2188.6
$ws.Range("I136").Value = 2232.111
$ws.Range("K136").Value = 6696.333
$ws.Range("M136").Value = -4146.333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1130.65
$ws.Range("I94").Value = 978.1429000000001
$ws.Range("K94").Value = 978.1429000000001
$ws.Range("M94").Value = -527.1429000000001
$ws.Range("H134").Value = 3979.5
$ws.Range("I134").Value = 3979.5
$ws.Range("K134").Value = 11938.5
$ws.Range("M134").Value = -9403.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 63533
$ws.Range("J108").Value = 63533
$ws.Range("L108").Value = 63533
$ws.Range("N108").Value = -71213
$ws.Range("H134").Value = 2387.5
$ws.Range("I134").Value = 550
$ws.Range("K134").Value = 1650
$ws.Range("M134").Value = 885
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1579.4
$ws.Range("I5").Value = 1579.4
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4738.200000000001
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4626.200000000001
$ws.Range("N5").Value = ""
$ws.Range("H116").Value = 999.5
$ws.Range("I116").Value = 999
$ws.Range("K116").Value = 2997
$ws.Range("M116").Value = 445
$ws.Range("H133").Value = 35207.625
$ws.Range("I133").Value = 124995
$ws.Range("J133").Value = 22380.857
$ws.Range("K133").Value = 374985
$ws.Range("L133").Value = 67142.571
$ws.Range("N133").Value = -77262.571
$ws.Range("M133").Value = -369925
$ws.Range("H134").Value = 11531.5
$ws.Range("I134").Value = 12635
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 37905
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -32835
$ws.Range("N134").Value = -14940
$ws.Range("H135").Value = 1579.4
$ws.Range("I135").Value = 1579.4
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 14214.6
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11679.6
$ws.Range("N135").Value = ""
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 44000
$ws.Range("J40").Value = 44000
$ws.Range("L40").Value = 44000
$ws.Range("N40").Value = -44302
$ws.Range("H102").Value = 2060.0908
$ws.Range("I102").Value = 2079.111
$ws.Range("J102").Value = 1974.5
$ws.Range("K102").Value = 2079.111
$ws.Range("L102").Value = 1974.5
$ws.Range("M102").Value = -457.1109999999999
$ws.Range("N102").Value = -5218.5
$ws.Range("H126").Value = 8428.5
$ws.Range("J126").Value = 7407
$ws.Range("L126").Value = 22221
$ws.Range("N126").Value = -27161
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2391.5454
$ws.Range("I22").Value = 2230
$ws.Range("J22").Value = 2674.25
$ws.Range("K22").Value = 2230
$ws.Range("L22").Value = 2674.25
$ws.Range("M22").Value = -1935
$ws.Range("N22").Value = -3264.25
$ws.Range("H27").Value = 2391.5454
$ws.Range("I27").Value = 2230
$ws.Range("J27").Value = 2674.25
$ws.Range("K27").Value = 2230
$ws.Range("L27").Value = 2674.25
$ws.Range("M27").Value = -2123
$ws.Range("N27").Value = -2888.25
$ws.Range("H40").Value = 3787
$ws.Range("I40").Value = 3078.2856
$ws.Range("K40").Value = 3078.2856
$ws.Range("M40").Value = -2942.2856
$ws.Range("H61").Value = 1125.2858
$ws.Range("I61").Value = 961
$ws.Range("J61").Value = 1536
$ws.Range("K61").Value = 961
$ws.Range("L61").Value = 1536
$ws.Range("M61").Value = -759
$ws.Range("N61").Value = -1940
$ws.Range("H82").Value = 1258.2
$ws.Range("I82").Value = 1126.2858
$ws.Range("K82").Value = 1126.2858
$ws.Range("M82").Value = -765.2858000000001
$ws.Range("H85").Value = 1258.2
$ws.Range("I85").Value = 1126.2858
$ws.Range("K85").Value = 1126.2858
$ws.Range("M85").Value = 121.7141999999999
$ws.Range("H93").Value = 758.1818
$ws.Range("I93").Value = 785
$ws.Range("J93").Value = 490
$ws.Range("K93").Value = 785
$ws.Range("L93").Value = 490
$ws.Range("M93").Value = 463
$ws.Range("N93").Value = -2986
$ws.Range("H100").Value = 1539.3334
$ws.Range("I100").Value = 1613
$ws.Range("K100").Value = 1613
$ws.Range("M100").Value = -1072
$ws.Range("H113").Value = 1125.2858
$ws.Range("I113").Value = 961
$ws.Range("J113").Value = 1536
$ws.Range("K113").Value = 961
$ws.Range("L113").Value = 1536
$ws.Range("M113").Value = 1209
$ws.Range("N113").Value = -5876
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value = 750
$ws.Range("J81").Value = 6667001
$ws.Range("K81").Value = 1500
$ws.Range("L81").Value = 13334002
$ws.Range("M81").Value = -439
$ws.Range("N81").Value = -13336124
$ws.Range("I84").Value = 750
$ws.Range("J84").Value = 6667001
$ws.Range("K84").Value = 7500
$ws.Range("L84").Value = 66670010
$ws.Range("M84").Value = -2196
$ws.Range("N84").Value = -66680618
$ws.Range("H100").Value = 12501622
$ws.Range("I100").Value = 14287471
$ws.Range("K100").Value = 28574942
$ws.Range("M100").Value = -28574401
$ws.Range("H107").Value = 1516.3125
$ws.Range("I107").Value = 1483.4546
$ws.Range("K107").Value = 4450.3638
$ws.Range("M107").Value = -2530.3638
$ws.Range("H122").Value = 6805.143
$ws.Range("I122").Value = 5991.3335
$ws.Range("K122").Value = 17974.0005
$ws.Range("M122").Value = -15524.0005
$ws.Range("H126").Value = 1375.0714
$ws.Range("I126").Value = 1271.3334
$ws.Range("K126").Value = 3814.0002
$ws.Range("M126").Value = -1344.0002
